$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Nome Messaggio"
$ws.Range("B1").Value = "Codice"
$ws.Range("C1").Value = "Descrizione"
$ws.Range("D1").Value = "Tipo*"
$ws.Range("F1").Value = "*"

# Data rows
$ws.Range("A2").Value = "ConnectionConfirmation"
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = "Messaggio di prova dal client per verificare la connessione"
$ws.Range("D2").Value = 2
$ws.Range("F2").Value = "0=>MultiDirezione"

$ws.Range("A3").Value = "ConnectionConfirmed"
$ws.Range("B3").Value = 1001
$ws.Range("C3").Value = "Messaggio dell'host per confermare il messaggio del client"
$ws.Range("D3").Value = 1
$ws.Range("F3").Value = "1=>Host verso client"

$ws.Range("A4").Value = "LobbyLogin"
$ws.Range("B4").Value = 1002
$ws.Range("C4").Value = "Messaggio per informare che il client è entrato nella lobby e deve ricevere le informazioni della lobby."
$ws.Range("D4").Value = 2
$ws.Range("F4").Value = "2=>Client verso Host"

$ws.Range("A5").Value = "LobbyInfo"
$ws.Range("B5").Value = 1003
$ws.Range("C5").Value = "Messaggio contenent le informazioni della lobby"
$ws.Range("D5").Value = 1

$ws.Range("A6").Value = "LobbyChatMessage"
$ws.Range("B6").Value = 1004
$ws.Range("C6").Value = "Manda un messggio di testo nella chat della lobby"
$ws.Range("D6").Value = 0

# Formatting - header row bold, size 16, height 21
$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("A1:F1").Font.Size = 16
$ws.Rows.Item(1).RowHeight = 21

# Underlined empty cell at C23
$ws.Range("C23").Font.Underline = $true

# Column widths
$ws.Columns.Item(1).ColumnWidth = 28.36328125
$ws.Columns.Item(2).ColumnWidth = 11.81640625
$ws.Columns.Item(3).ColumnWidth = 113.7265625
$ws.Columns.Item(4).ColumnWidth = 17
$ws.Columns.Item(6).ColumnWidth = 24.7265625

# Selection
$ws.Range("C6").Select()
